# Auto-generated Excel COM-interop edit script
# Applies numeric cell-value corrections to the Kraken_Profits workbook
# (sheets ALC, ARM, BSM, CRP, CUL, LTW, WVR) per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 359.625
$ws.Range("I2").Value = 359.625
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 359.625
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -246.625
$ws.Range("N2").ClearContents()
$ws.Range("H38").Value = 11034.667
$ws.Range("I38").Value = 1388.2727
$ws.Range("K38").Value = 4164.8181
$ws.Range("M38").Value = -3792.8181
$ws.Range("H43").Value = 4000.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 4000.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 4000.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -4138.5
$ws.Range("H135").Value = 1578.8334
$ws.Range("I135").Value = 1578.8334
$ws.Range("K135").Value = 14209.5006
$ws.Range("M135").Value = -11674.5006

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 655.5
$ws.Range("I2").Value = 655.5
$ws.Range("K2").Value = 655.5
$ws.Range("M2").Value = -542.5
$ws.Range("H61").Value = 3100
$ws.Range("I61").Value = 3100
$ws.Range("K61").Value = 3100
$ws.Range("M61").Value = -2888
$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14126
$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 75000
$ws.Range("M77").Value = -70632
$ws.Range("H110").Value = 4666
$ws.Range("I110").Value = 4499.5
$ws.Range("K110").Value = 4499.5
$ws.Range("M110").Value = -2454.5
$ws.Range("H116").Value = 655.5
$ws.Range("I116").Value = 655.5
$ws.Range("K116").Value = 655.5
$ws.Range("M116").Value = 1638.5
$ws.Range("H136").Value = 3100
$ws.Range("I136").Value = 3100
$ws.Range("K136").Value = 9300
$ws.Range("M136").Value = -6750

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 655.5
$ws.Range("I3").Value = 655.5
$ws.Range("K3").Value = 655.5
$ws.Range("M3").Value = -541.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H122").Value = 966.6667
$ws.Range("I122").Value = 966.6667
$ws.Range("K122").Value = 2900.0001
$ws.Range("M122").Value = -450.0001000000002
$ws.Range("H132").Value = 4689.125
$ws.Range("I132").Value = 3928.4285
$ws.Range("J132").Value = 10014
$ws.Range("K132").Value = 11785.2855
$ws.Range("L132").Value = 30042
$ws.Range("M132").Value = -9255.2855
$ws.Range("N132").Value = -35102

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 58.583332
$ws.Range("I2").Value = 47.625
$ws.Range("J2").Value = 80.5
$ws.Range("K2").Value = 285.75
$ws.Range("L2").Value = 483
$ws.Range("M2").Value = -172.75
$ws.Range("N2").Value = -709
$ws.Range("H14").Value = 900
$ws.Range("I14").Value = 900
$ws.Range("K14").Value = 2700
$ws.Range("M14").Value = -2527
$ws.Range("H23").Value = 710.6
$ws.Range("I23").Value = 542.9167
$ws.Range("J23").Value = 865.38464
$ws.Range("K23").Value = 1628.7501
$ws.Range("L23").Value = 2596.15392
$ws.Range("M23").Value = -1393.7501
$ws.Range("N23").Value = -3066.15392
$ws.Range("H33").Value = 917.6
$ws.Range("I33").Value = 949.5
$ws.Range("J33").Value = 896.3333
$ws.Range("K33").Value = 5697
$ws.Range("L33").Value = 5377.9998
$ws.Range("M33").Value = -5414
$ws.Range("N33").Value = -5943.9998
$ws.Range("H38").Value = 48.6
$ws.Range("I38").Value = 91
$ws.Range("J38").Value = 38
$ws.Range("K38").Value = 273
$ws.Range("L38").Value = 114
$ws.Range("M38").Value = 74
$ws.Range("N38").Value = -808
$ws.Range("H86").Value = 14250.5
$ws.Range("I86").Value = 1000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1814
$ws.Range("H89").Value = 14250.5
$ws.Range("I89").Value = 1000
$ws.Range("K89").Value = 9000
$ws.Range("M89").Value = -3072
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H129").Value = 3135.1333
$ws.Range("I129").Value = 2285.1428
$ws.Range("K129").Value = 6855.428400000001
$ws.Range("M129").Value = -1855.428400000001
$ws.Range("H131").Value = 1355.7142
$ws.Range("I131").Value = 1098
$ws.Range("K131").Value = 3294
$ws.Range("M131").Value = 1746
$ws.Range("H136").Value = 1754.6
$ws.Range("I136").Value = 1754.6
$ws.Range("K136").Value = 5263.799999999999
$ws.Range("M136").Value = -163.7999999999993
$ws.Range("H138").Value = 2204.2
$ws.Range("I138").Value = 2204.2
$ws.Range("K138").Value = 6612.599999999999
$ws.Range("M138").Value = -1472.599999999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2814.2778
$ws.Range("J22").Value = 3150.6
$ws.Range("L22").Value = 3150.6
$ws.Range("N22").Value = -3740.6
$ws.Range("H27").Value = 2814.2778
$ws.Range("J27").Value = 3150.6
$ws.Range("L27").Value = 3150.6
$ws.Range("N27").Value = -3364.6
$ws.Range("H40").Value = 17169.166
$ws.Range("I40").Value = 16603
$ws.Range("J40").Value = 20000
$ws.Range("K40").Value = 16603
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = -16467
$ws.Range("N40").Value = -20272
$ws.Range("H122").Value = 7333
$ws.Range("I122").Value = 8000
$ws.Range("J122").Value = 5999
$ws.Range("K122").Value = 24000
$ws.Range("L122").Value = 17997
$ws.Range("M122").Value = -21550
$ws.Range("N122").Value = -22897

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
